# Auto-generated Excel COM-interop edit script
# Applies scheduled-runner profit recalculation updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 646.13336
$ws.Range("I19").Value = 523.63635
$ws.Range("K19").Value = 523.63635
$ws.Range("M19").Value = -348.63635

$ws.Range("H51").Value = 5020.1055
$ws.Range("I51").Value = 4750
$ws.Range("J51").Value = 5051.8823
$ws.Range("K51").Value = 4750
$ws.Range("L51").Value = 5051.8823
$ws.Range("M51").Value = -4266
$ws.Range("N51").Value = -6019.8823

$ws.Range("H121").Value = 773.95654
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 790.9545000000001
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 2372.8635
$ws.Range("M121").Value = 547
$ws.Range("N121").Value = -5866.8635

$ws.Range("H129").Value = 1426.5454
$ws.Range("J129").Value = 1989.1428
$ws.Range("L129").Value = 5967.428400000001
$ws.Range("N129").Value = -15967.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 113090
$ws.Range("I2").Value = 127101.25
$ws.Range("K2").Value = 127101.25
$ws.Range("M2").Value = -126988.25

$ws.Range("H45").Value = 970.6667
$ws.Range("I45").Value = 956
$ws.Range("K45").Value = 956
$ws.Range("M45").Value = -579

$ws.Range("H61").Value = 2242.9644
$ws.Range("I61").Value = 1556.238
$ws.Range("J61").Value = 4303.143
$ws.Range("K61").Value = 1556.238
$ws.Range("L61").Value = 4303.143
$ws.Range("M61").Value = -1344.238
$ws.Range("N61").Value = -4727.143

$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

$ws.Range("H116").Value = 113090
$ws.Range("I116").Value = 127101.25
$ws.Range("K116").Value = 127101.25
$ws.Range("M116").Value = -124807.25

$ws.Range("H132").Value = 2612.1091
$ws.Range("I132").Value = 2026.2727
$ws.Range("J132").Value = 4955.4546
$ws.Range("K132").Value = 6078.8181
$ws.Range("L132").Value = 14866.3638
$ws.Range("M132").Value = -3548.8181
$ws.Range("N132").Value = -19926.3638

$ws.Range("H136").Value = 2242.9644
$ws.Range("I136").Value = 1556.238
$ws.Range("J136").Value = 4303.143
$ws.Range("K136").Value = 4668.714
$ws.Range("L136").Value = 12909.429
$ws.Range("M136").Value = -2118.714
$ws.Range("N136").Value = -18009.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 113090
$ws.Range("I3").Value = 127101.25
$ws.Range("K3").Value = 127101.25
$ws.Range("M3").Value = -126987.25

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3228.395
$ws.Range("I31").Value = 1448.825
$ws.Range("J31").Value = 4964.561
$ws.Range("K31").Value = 1448.825
$ws.Range("L31").Value = 4964.561
$ws.Range("M31").Value = -1153.825
$ws.Range("N31").Value = -5554.561

$ws.Range("H34").Value = 3228.395
$ws.Range("I34").Value = 1448.825
$ws.Range("J34").Value = 4964.561
$ws.Range("K34").Value = 1448.825
$ws.Range("L34").Value = 4964.561
$ws.Range("M34").Value = -1246.825
$ws.Range("N34").Value = -5368.561

$ws.Range("H58").Value = 33335154
$ws.Range("I58").Value = 43479424
$ws.Range("J58").Value = 3973
$ws.Range("K58").Value = 43479424
$ws.Range("L58").Value = 3973
$ws.Range("M58").Value = -43479221
$ws.Range("N58").Value = -4379

$ws.Range("H99").Value = 18520966
$ws.Range("I99").Value = 1881.8889
$ws.Range("J99").Value = 37040052
$ws.Range("K99").Value = 1881.8889
$ws.Range("L99").Value = 37040052
$ws.Range("M99").Value = -383.8888999999999
$ws.Range("N99").Value = -37043048

$ws.Range("H105").Value = 801.16327
$ws.Range("I105").Value = 759.6579
$ws.Range("K105").Value = 759.6579
$ws.Range("M105").Value = 987.3421

$ws.Range("H126").Value = 18520966
$ws.Range("I126").Value = 1881.8889
$ws.Range("J126").Value = 37040052
$ws.Range("K126").Value = 5645.6667
$ws.Range("L126").Value = 111120156
$ws.Range("M126").Value = -3175.6667
$ws.Range("N126").Value = -111125096

$ws.Range("H132").Value = 8335945.5
$ws.Range("I132").Value = 13890480
$ws.Range("K132").Value = 41671440
$ws.Range("M132").Value = -41668910

$ws.Range("H134").Value = 33335684
$ws.Range("I134").Value = 35715720
$ws.Range("J134").Value = 29415624
$ws.Range("K134").Value = 107147160
$ws.Range("L134").Value = 88246872
$ws.Range("M134").Value = -107144625
$ws.Range("N134").Value = -88251942

$ws.Range("H136").Value = 33335154
$ws.Range("I136").Value = 43479424
$ws.Range("J136").Value = 3973
$ws.Range("K136").Value = 130438272
$ws.Range("L136").Value = 11919
$ws.Range("M136").Value = -130435722
$ws.Range("N136").Value = -17019

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 11087.714
$ws.Range("I68").Value = 11272.25
$ws.Range("J68").Value = 9980.5
$ws.Range("K68").Value = 33816.75
$ws.Range("L68").Value = 29941.5
$ws.Range("M68").Value = -33005.75
$ws.Range("N68").Value = -31563.5

$ws.Range("H71").Value = 11087.714
$ws.Range("I71").Value = 11272.25
$ws.Range("J71").Value = 9980.5
$ws.Range("K71").Value = 101450.25
$ws.Range("L71").Value = 89824.5
$ws.Range("M71").Value = -97394.25
$ws.Range("N71").Value = -97936.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 65.333336
$ws.Range("I2").Value = 74
$ws.Range("K2").Value = 74
$ws.Range("M2").Value = 39

$ws.Range("H138").Value = 58857.25
$ws.Range("J138").Value = 58857.25
$ws.Range("L138").Value = 58857.25
$ws.Range("N138").Value = -69137.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 250.54546
$ws.Range("I55").Value = 239.17647
$ws.Range("J55").Value = 289.2
$ws.Range("K55").Value = 239.17647
$ws.Range("L55").Value = 289.2
$ws.Range("M55").Value = -66.17646999999999
$ws.Range("N55").Value = -635.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 262582.25
$ws.Range("J46").Value = 262582.25
$ws.Range("L46").Value = 262582.25
$ws.Range("N46").Value = -263044.25

$ws.Range("H121").Value = 19990
$ws.Range("J121").Value = 19990
$ws.Range("L121").Value = 19990
$ws.Range("N121").Value = -23484

$ws.Range("H134").Value = 262582.25
$ws.Range("J134").Value = 262582.25
$ws.Range("L134").Value = 787746.75
$ws.Range("N134").Value = -792816.75
